# Fruta / hortaliza, semanal
# Insert a new weekly record at row 24, shifting the existing rows 24-52
# down to rows 25-53. The new row reuses the constant (market/category)
# columns from the row that is pushed down to row 25, and carries its own
# date (D) and volume (J) figures, with price columns (K/L/M/P) matching
# the existing "12000 / 1200" price band.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 24..52 down by one row, creating a blank row 24.
$ws.Rows.Item(24).Insert()

# Fill the new row 24. Non numeric/date columns are copied from the row
# that used to occupy row 24 (now shifted to row 25), since those columns
# (market, region, category, etc.) are constant throughout the sheet.
$ws.Range("A24").Value = $ws.Range("A25").Value()
$ws.Range("B24").Value = $ws.Range("B25").Value()
$ws.Range("C24").Value = $ws.Range("C25").Value()
$ws.Range("D24").Value = 44848
$ws.Range("E24").Value = $ws.Range("E25").Value()
$ws.Range("F24").Value = $ws.Range("F25").Value()
$ws.Range("G24").Value = $ws.Range("G25").Value()
$ws.Range("H24").Value = $ws.Range("H25").Value()
$ws.Range("I24").Value = $ws.Range("I25").Value()
$ws.Range("J24").Value = 35
$ws.Range("K24").Value = 12000
$ws.Range("L24").Value = 12000
$ws.Range("M24").Value = 12000
$ws.Range("N24").Value = $ws.Range("N25").Value()
$ws.Range("O24").Value = $ws.Range("O25").Value()
$ws.Range("P24").Value = 1200
$ws.Range("Q24").Value = $ws.Range("Q25").Value()
$ws.Range("R24").Value = $ws.Range("R25").Value()
